$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "n_sikap_a"
$ws.Range("C1").Value = "mother_work_lainnya"
$ws.Range("H1").Value = "mother_salary_sangat_rendah"
$ws.Range("I1").Value = "father_salary_tidak_berpenghasilan"
$ws.Range("L1").Value = "father_edu_smp_sederajat"
$ws.Range("M1").Value = "father_work_buruh"
$ws.Range("N1").Value = "mother_salary_cukup_rendah"
$ws.Range("O1").Value = "mother_work_buruh"
